$wb = $excel.ActiveWorkbook

# Totales Plantel 2P
$ws2 = $wb.Worksheets.Item("Totales Plantel 2P")
$ws2.Range("G3").Value = 1
$ws2.Range("H3").Value = 2
$ws2.Range("I3").Value = 1
$ws2.Range("J3").Value = 1
$ws2.Range("K3").Value = 4
$ws2.Range("L3").Value = 11
$ws2.Range("M3").Value = 14
$ws2.Range("N3").Value = 157
$ws2.Range("O3").Value = 34
$ws2.Range("P3").Value = 17.8
$ws2.Range("J4").Value = 1
$ws2.Range("K4").Value = 10
$ws2.Range("L4").Value = 9
$ws2.Range("M4").Value = 16
$ws2.Range("N4").Value = 173
$ws2.Range("O4").Value = 37
$ws2.Range("P4").Value = 17.6
$ws2.Range("G5").Value = 5
$ws2.Range("H5").Value = 4
$ws2.Range("I5").Value = 10
$ws2.Range("J5").Value = 11
$ws2.Range("K5").Value = 22
$ws2.Range("L5").Value = 33
$ws2.Range("M5").Value = 83
$ws2.Range("N5").Value = 465
$ws2.Range("O5").Value = 170
$ws2.Range("P5").Value = 26.8
$ws2.Range("G6").Value = 2
$ws2.Range("I6").Value = 4
$ws2.Range("J6").Value = 12
$ws2.Range("K6").Value = 19
$ws2.Range("L6").Value = 13
$ws2.Range("M6").Value = 12
$ws2.Range("N6").Value = 62
$ws2.Range("O6").Value = 66
$ws2.Range("P6").Value = 51.6
$ws2.Range("K7").Value = 6
$ws2.Range("L7").Value = 11
$ws2.Range("M7").Value = 19
$ws2.Range("N7").Value = 42
$ws2.Range("O7").Value = 41
$ws2.Range("P7").Value = 49.4
$ws2.Range("I8").Value = 0
$ws2.Range("J8").Value = 4
$ws2.Range("L8").Value = 5
$ws2.Range("M8").Value = 9
$ws2.Range("N8").Value = 89
$ws2.Range("O8").Value = 26
$ws2.Range("P8").Value = 22.6
$ws2.Range("G9").Value = 4
$ws2.Range("I9").Value = 5
$ws2.Range("J9").Value = 18
$ws2.Range("K9").Value = 33
$ws2.Range("L9").Value = 29
$ws2.Range("M9").Value = 40
$ws2.Range("N9").Value = 193
$ws2.Range("O9").Value = 133
$ws2.Range("P9").Value = 40.8
$ws2.Range("G10").Value = 9
$ws2.Range("H10").Value = 8
$ws2.Range("I10").Value = 15
$ws2.Range("J10").Value = 29
$ws2.Range("K10").Value = 55
$ws2.Range("L10").Value = 62
$ws2.Range("M10").Value = 123
$ws2.Range("N10").Value = 658
$ws2.Range("O10").Value = 303
$ws2.Range("P10").Value = 31.5

# Totales Plantel Final
$ws3 = $wb.Worksheets.Item("Totales Plantel Final")
$ws3.Range("I3").Value = 1
$ws3.Range("J3").Value = 1
$ws3.Range("K3").Value = 4
$ws3.Range("L3").Value = 11
$ws3.Range("M3").Value = 14
$ws3.Range("N3").Value = 157
$ws3.Range("O3").Value = 34
$ws3.Range("P3").Value = 17.8
$ws3.Range("J4").Value = 1
$ws3.Range("K4").Value = 10
$ws3.Range("I5").Value = 8
$ws3.Range("J5").Value = 11
$ws3.Range("K5").Value = 19
$ws3.Range("L5").Value = 31
$ws3.Range("M5").Value = 58
$ws3.Range("N5").Value = 497
$ws3.Range("O5").Value = 138
$ws3.Range("P5").Value = 21.7
$ws3.Range("G6").Value = 2
$ws3.Range("H6").Value = 4
$ws3.Range("I6").Value = 4
$ws3.Range("J6").Value = 12
$ws3.Range("K6").Value = 19
$ws3.Range("L6").Value = 13
$ws3.Range("M6").Value = 12
$ws3.Range("N6").Value = 62
$ws3.Range("O6").Value = 66
$ws3.Range("P6").Value = 51.6
$ws3.Range("M8").Value = 9
$ws3.Range("N8").Value = 89
$ws3.Range("O8").Value = 26
$ws3.Range("P8").Value = 22.6
$ws3.Range("G9").Value = 4
$ws3.Range("H9").Value = 4
$ws3.Range("I9").Value = 5
$ws3.Range("J9").Value = 18
$ws3.Range("K9").Value = 33
$ws3.Range("L9").Value = 29
$ws3.Range("M9").Value = 40
$ws3.Range("N9").Value = 193
$ws3.Range("O9").Value = 133
$ws3.Range("P9").Value = 40.8
$ws3.Range("G10").Value = 9
$ws3.Range("H10").Value = 8
$ws3.Range("I10").Value = 13
$ws3.Range("J10").Value = 29
$ws3.Range("K10").Value = 52
$ws3.Range("L10").Value = 60
$ws3.Range("M10").Value = 98
$ws3.Range("N10").Value = 690
$ws3.Range("O10").Value = 271
$ws3.Range("P10").Value = 28.2

# Reprobados por Grupo
$ws4 = $wb.Worksheets.Item("Reprobados por Grupo")
$ws4.Range("J8").Value = 1
$ws4.Range("L8").Value = 3
$ws4.Range("M8").Value = 15
$ws4.Range("N8").Value = 9
$ws4.Range("O8").Value = 37.5
$ws4.Range("H10").Value = 0
$ws4.Range("I10").Value = 1
$ws4.Range("K12").Value = 6
$ws4.Range("L12").Value = 3
$ws4.Range("K13").Value = 1
$ws4.Range("L13").Value = 3
$ws4.Range("I17").Value = 0
$ws4.Range("J17").Value = 2
$ws4.Range("J20").Value = 8
$ws4.Range("K20").Value = 5
$ws4.Range("L20").Value = 2
$ws4.Range("M20").Value = 13
$ws4.Range("N20").Value = 23
$ws4.Range("O20").Value = 63.9
$ws4.Range("H21").Value = 0
$ws4.Range("I21").Value = 2
$ws4.Range("K21").Value = 5
$ws4.Range("L21").Value = 2
$ws4.Range("M21").Value = 15
$ws4.Range("N21").Value = 13
$ws4.Range("O21").Value = 46.4
$ws4.Range("F22").Value = 2
$ws4.Range("G22").Value = 2
$ws4.Range("H22").Value = 0
$ws4.Range("M22").Value = 10
$ws4.Range("N22").Value = 13
$ws4.Range("O22").Value = 56.5
$ws4.Range("H23").Value = 1
$ws4.Range("I23").Value = 3
$ws4.Range("J23").Value = 5
$ws4.Range("M23").Value = 17
$ws4.Range("N23").Value = 13
$ws4.Range("O23").Value = 43.3
$ws4.Range("L33").Value = 1
$ws4.Range("M33").Value = 24
$ws4.Range("N33").Value = 1
$ws4.Range("O33").Value = 4

# Totales Grupos
$ws5 = $wb.Worksheets.Item("Totales Grupos")
$ws5.Range("D8").Value = 15
$ws5.Range("E8").Value = 15
$ws5.Range("G8").Value = 62.5
$ws5.Range("H8").Value = 62.5
$ws5.Range("D9").Value = 30
$ws5.Range("G9").Value = 88.2
$ws5.Range("D10").Value = 27
$ws5.Range("G10").Value = 87.09999999999999
$ws5.Range("D11").Value = 38
$ws5.Range("G11").Value = 95
$ws5.Range("D12").Value = 14
$ws5.Range("G12").Value = 58.3
$ws5.Range("D13").Value = 33
$ws5.Range("G13").Value = 86.8
$ws5.Range("D17").Value = 32
$ws5.Range("G17").Value = 82.09999999999999
$ws5.Range("D20").Value = 13
$ws5.Range("E20").Value = 13
$ws5.Range("G20").Value = 36.1
$ws5.Range("H20").Value = 36.1
$ws5.Range("D21").Value = 15
$ws5.Range("E21").Value = 15
$ws5.Range("G21").Value = 53.6
$ws5.Range("H21").Value = 53.6
$ws5.Range("D22").Value = 10
$ws5.Range("E22").Value = 10
$ws5.Range("G22").Value = 43.5
$ws5.Range("H22").Value = 43.5
$ws5.Range("D23").Value = 17
$ws5.Range("E23").Value = 17
$ws5.Range("G23").Value = 56.7
$ws5.Range("H23").Value = 56.7
$ws5.Range("D24").Value = 7
$ws5.Range("G24").Value = 63.6
$ws5.Range("D28").Value = 12
$ws5.Range("G28").Value = 80
$ws5.Range("D29").Value = 7
$ws5.Range("G29").Value = 63.6
$ws5.Range("D32").Value = 14
$ws5.Range("G32").Value = 82.40000000000001
$ws5.Range("D33").Value = 24
$ws5.Range("E33").Value = 24
$ws5.Range("G33").Value = 96
$ws5.Range("H33").Value = 96
$ws5.Range("D34").Value = 23
$ws5.Range("G34").Value = 92
